$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81, shifting existing rows 81-106 down to 82-107.
$ws.Rows.Item(81).Insert()

# Populate the newly inserted row 81 with the new weekly data point.
$ws.Cells.Item(81, 1).Value = 1
$ws.Cells.Item(81, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(81, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(81, 4).Value = 44841
$ws.Cells.Item(81, 5).Value = 15
$ws.Cells.Item(81, 6).Value = 100112038
$ws.Cells.Item(81, 7).Value = "Cebollín baby"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 300
$ws.Cells.Item(81, 11).Value = 900
$ws.Cells.Item(81, 12).Value = 1000
$ws.Cells.Item(81, 13).Value = 950
$ws.Cells.Item(81, 14).Value = "`$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(81, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(81, 16).Value = 475
$ws.Cells.Item(81, 17).Value = 2
$ws.Cells.Item(81, 18).Value = "Hortaliza"
